$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / percentage / name / URL updates (safe from numeric auto-conversion)
$ws.Range("D2").Value = "64.226.98"
$ws.Range("E2").Value = "  -3.34%  "
$ws.Range("D3").Value = "3.178.07"
$ws.Range("E3").Value = "  -4.12%  "
$ws.Range("E5").Value = "  -3.06%  "
$ws.Range("E6").Value = "  -7.86%  "
$ws.Range("E7").Value = "  -5.47%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "3.177.29"
$ws.Range("E9").Value = "  -4.07%  "
$ws.Range("E10").Value = "  -4.56%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("E12").Value = "  -5.00%  "
$ws.Range("D13").Value = "3.730.85"
$ws.Range("E13").Value = "  -3.95%  "
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "64.308.31"
$ws.Range("E15").Value = "  -3.23%  "
$ws.Range("E16").Value = "  -4.10%  "
$ws.Range("E17").Value = "  -3.86%  "
$ws.Range("D18").Value = "3.176.10"
$ws.Range("E18").Value = "  -4.12%  "
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E20").Value = "  -3.48%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E21").Value = "  -4.26%  "
$ws.Range("E22").Value = "  -4.95%  "
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("E24").Value = "  -2.86%  "
$ws.Range("E25").Value = "  +2.74%  "
$ws.Range("E26").Value = "  -4.98%  "
$ws.Range("E27").Value = "  -7.28%  "
$ws.Range("E28").Value = "  -3.12%  "
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E30").Value = "  -2.83%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E31").Value = "  -6.96%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E33").Value = "  -4.66%  "
$ws.Range("E34").Value = "  -4.59%  "
$ws.Range("E35").Value = "  -6.15%  "
$ws.Range("E36").Value = "  -2.24%  "
$ws.Range("E37").Value = "  -6.56%  "
$ws.Range("D38").Value = "2.697.93"
$ws.Range("E38").Value = "  -5.66%  "
$ws.Range("E39").Value = "  -7.97%  "
$ws.Range("E40").Value = "  -8.98%  "
$ws.Range("E41").Value = "  -4.43%  "
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("E43").Value = "  -8.51%  "
$ws.Range("E44").Value = "  -5.72%  "
$ws.Range("E45").Value = "  -6.72%  "
$ws.Range("E46").Value = "  -7.25%  "
$ws.Range("E47").Value = "  -7.32%  "
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("E49").Value = "  -12.82%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E51").Value = "  -4.81%  "

# Numeric-looking text values: force Text format so Excel does not coerce them to numbers,
# then restore the default "Normal" style so no stray formatting is introduced.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.606"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.120"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "418.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.490"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000106"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.987"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "155.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.708"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0623"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "295.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0261"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0991"
$ws.Range("D51").Style = "Normal"
